# Auto-update draw results: append the 2025-12-23 Pick 3 draw as a new row.
#
# The source workbook stores every data cell as literal text (dates like
# "2025-12-22" and numeric-looking values like "251222" are all plain
# strings, never real numbers/dates). To keep the new row consistent with
# the existing 97 rows, force the new cells to Text format before writing
# the values so Excel does not reinterpret "2025-12-23" as a date or
# "251223" as a number. The format is then reset back to the sheet's
# normal (General) style so the new row doesn't end up visually
# different from the rows above it - only the underlying text values
# need to be preserved verbatim.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 98
$rangeAddress = "A" + $newRow + ":E" + $newRow
$rowRange = $ws.Range($rangeAddress)

# Force text storage so the date- and number-looking strings are kept as-is
# (otherwise Excel would reinterpret "2025-12-23" as a date and "251223" as
# a number).
$rowRange.NumberFormat = "@"

$ws.Range("A" + $newRow).Value = "2025-12-23"
$ws.Range("B" + $newRow).Value = "Pick 3"
$ws.Range("C" + $newRow).Value = "251223"
$ws.Range("D" + $newRow).Value = "9-9-8"
$ws.Range("E" + $newRow).Value = "2025-12-23T21:42:44.985+04:00"

# Restore the default/normal style so the appended row matches the
# formatting of the rest of the table (only the text value itself needs
# to stay literal, not the cell format).
$rowRange.Style = "Normal"
